$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) store numeric-looking values as text
# (t="inlineStr" in the original). Temporarily format the range as Text so
# Excel does not auto-convert assigned strings into numbers, then restore
# the default "Normal" style so no stray number-format style is left behind.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.029.01'
$ws.Range("E2").Value = '  -0.55%  '

$ws.Range("D3").Value = '1.831.99'
$ws.Range("E3").Value = '  -0.42%  '

$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '242.16'
$ws.Range("E5").Value = '  -0.14%  '

$ws.Range("D6").Value = '0.6250'
$ws.Range("E6").Value = '  -5.68%  '

$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").Value = '0.07568'
$ws.Range("E8").Value = '  +1.69%  '

$ws.Range("D9").Value = '0.2918'
$ws.Range("E9").Value = '  -1.31%  '

$ws.Range("D10").Value = '22.53'
$ws.Range("E10").Value = '  -3.01%  '

$ws.Range("D11").Value = '0.07715'
$ws.Range("E11").Value = '  -0.35%  '

$ws.Range("D12").Value = '1.833.87'
$ws.Range("E12").Value = '  -0.35%  '

$ws.Range("D13").Value = '4.948'
$ws.Range("E13").Value = '  -1.31%  '

$ws.Range("D14").Value = '0.6636'
$ws.Range("E14").Value = '  -1.41%  '

$ws.Range("D15").Value = '0.00001013'
$ws.Range("E15").Value = '  +16.36%  '

$ws.Range("D16").Value = '82.55'
$ws.Range("E16").Value = '  -0.70%  '

$ws.Range("D17").Value = '6.021'
$ws.Range("E17").Value = '  -2.39%  '

$ws.Range("D18").Value = '29.021.08'
$ws.Range("E18").Value = '  -0.64%  '

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '226.14'
$ws.Range("E19").Value = '  -0.11%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '12.33'
$ws.Range("E20").Value = '  -1.52%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '0.9992'
$ws.Range("E21").Value = '  -0.19%  '

$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").Value = '7.155'
$ws.Range("E22").Value = '  -0.25%  '

$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").Value = '158.13'
$ws.Range("E24").Value = '  -0.49%  '

$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '8.478'
$ws.Range("E25").Value = '  -1.56%  '

$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '0.1371'
$ws.Range("E26").Value = '  -1.54%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '17.90'
$ws.Range("E27").Value = '  -0.75%  '

$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = '1.488'
$ws.Range("E28").Value = '  -1.61%  '

$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").Value = '4.093'
$ws.Range("E29").Value = '  -1.11%  '

$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '4.011'
$ws.Range("E30").Value = '  -0.74%  '

$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '1.191'
$ws.Range("E31").Value = '  -1.39%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '0.05200'
$ws.Range("E32").Value = '  -3.44%  '

$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").Value = '1.841'
$ws.Range("E33").Value = '  -0.82%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '0.7360'
$ws.Range("E34").Value = '  -1.64%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.137'
$ws.Range("E35").Value = '  -1.79%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.696'
$ws.Range("E36").Value = '  +1.66%  '

$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.237.71'
$ws.Range("E37").Value = '  -4.86%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '2.760'
$ws.Range("E38").Value = '  -0.04%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01781'
$ws.Range("E39").Value = '  -0.80%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '6.325'
$ws.Range("E40").Value = '  -0.42%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '0.8949'
$ws.Range("E41").Value = '  -1.04%  '

$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '0.9995'
$ws.Range("E42").Value = '  -0.09%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '101.38'
$ws.Range("E43").Value = '  -2.32%  '

$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '1.978.53'
$ws.Range("E44").Value = '  -0.52%  '

$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.00000000123'
$ws.Range("E45").Value = '  -0.71%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '64.14'
$ws.Range("E46").Value = '  -1.56%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '0.5106'
$ws.Range("E47").Value = '  -0.69%  '

$ws.Range("B48").Value = 'TheSandbox'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D48").Value = '0.4024'
$ws.Range("E48").Value = '  +0.05%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '8.848'
$ws.Range("E49").Value = '  +1.17%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.05753'
$ws.Range("E50").Value = '  -1.83%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '1.638'
$ws.Range("E51").Value = '  -6.51%  '

$priceVolRange.Style = "Normal"

